$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the per-row "Id" (column A) GUIDs for existing rows 2-10 ---
$ws.Range("A2").Value = "4e4a6a37-9b97-4298-8ed4-35ddf1601bd1"
$ws.Range("A3").Value = "62b8bad1-1fcf-487a-b220-199e1070b10b"
$ws.Range("A4").Value = "8c4abf03-63a4-4879-bd99-9b6be0eea067"
$ws.Range("A5").Value = "e4a7f74e-bd07-4e51-ad52-cff71e98d353"
$ws.Range("A6").Value = "c8cc7c79-ae8b-4651-9054-d6093c81474e"
$ws.Range("A7").Value = "844a3dd0-a9b6-46b4-929f-177d9076fe27"
$ws.Range("A8").Value = "6418f889-1a37-411b-aceb-8cd7501d3c82"
$ws.Range("A9").Value = "758a38f5-59b0-49d5-9387-c0baa6d35f8d"
$ws.Range("A10").Value = "39b5a791-f34a-48d2-9eb7-a85b3e04b6c3"

# --- Update the "UserId" (column G) GUID, shared by all rows 2-10 ---
$ws.Range("G2:G10").Value = "c09cdc72-6c50-4ff4-9a9e-5fdbc855e8f3"

# --- Update numeric stats for rows 2, 3 and 7 ---
$ws.Range("D2").Value = 251
$ws.Range("E2").Value = 5036
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 45

$ws.Range("D3").Value = 333
$ws.Range("H3").Value = 336
$ws.Range("I3").Value = 326

$ws.Range("D7").Value = 376
$ws.Range("E7").Value = 657
$ws.Range("H7").Value = 195
$ws.Range("I7").Value = 187

# --- Add new row 11 for summoner "KnifeTheSkull" ---
$ws.Range("A11").Value = "3db97f18-977c-4821-bf28-b7e255c4b3e1"
$ws.Range("B11").Value = "KnifeTheSkull"
$ws.Range("C11").Value = "duryWWTOrWpp7wrWyiKPICFY59CvT1W5KfO329u3e6d6vQA"
$ws.Range("D11").Value = 227
$ws.Range("E11").Value = 26
$ws.Range("F11").Value = "SjwN75jquERgi0ch1qfBF06Y-dVtHF5HerErgWj6_2LVbw"
$ws.Range("G11").Value = "c09cdc72-6c50-4ff4-9a9e-5fdbc855e8f3"
$ws.Range("H11").Value = 90
$ws.Range("I11").Value = 98
